# Apply the new "grafics" class columns (E:G) which duplicate the
# existing B:D data, and rename the old B1/C1 headers to the ".0" variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers C1 ("l25" -> "l25.0") and B1 ("l23" -> "l23.0")
$ws.Range("C1").Value = "l25.0"
$ws.Range("B1").Value = "l23.0"

# Add the new headers for the duplicated columns E, F, G
$ws.Range("E1").Value = "l23"
$ws.Range("F1").Value = "l25"
$ws.Range("G1").Value = "l27"

# Copy the data values from B2:D12 into the new E2:G12 range
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 5).Value = $ws.Cells.Item($row, 2).Value2  # E = B
    $ws.Cells.Item($row, 6).Value = $ws.Cells.Item($row, 3).Value2  # F = C
    $ws.Cells.Item($row, 7).Value = $ws.Cells.Item($row, 4).Value2  # G = D
}

# Update the active selection to reflect the newly edited cell G2
$ws.Range("G2").Select()
